$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the Stock_date column (B2:B7) from a custom date-format value
# ("02-29-2022") to a plain text value ("29-Apr-2022").
$ws.Range("B2:B7").NumberFormat = "@"
$ws.Range("B2:B7").Value = "29-Apr-2022"

# Re-enter the amount formula across I2:I7 so Excel collapses it back into
# a shared-formula group.
$ws.Range("I2:I7").Formula = "=(F:F*H:H)"

# Update the active selection to B5.
[void]$ws.Range("B5").Select()
